# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Column A = rank index (unchanged), B = coin name, C = coinranking.com link,
# D = price, E = 1h volume change. Two coin rows (WrappedEther/Litecoin and
# Toncoin/Stellar) swapped rank position, so B/C/D/E are rewritten for those rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    # D/E cells store text that often LOOKS numeric ("243.82", "0.4709", "  -1.79%  ").
    # The source workbook keeps them as plain text (inline strings), so force the
    # Text format before assigning -- otherwise Excel auto-coerces the value to a
    # Number and mangles formatting (drops padding/zeros, switches to scientific, etc).
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell 2 4 '30.378.69'
Set-TextCell 2 5 '  -0.79%  '
Set-TextCell 3 4 '1.870.44'
Set-TextCell 3 5 '  -0.41%  '
Set-TextCell 4 5 '  -0.01%  '
Set-TextCell 5 4 '243.82'
Set-TextCell 5 5 '  -1.79%  '
Set-TextCell 6 5 '  -0.01%  '
Set-TextCell 7 4 '0.4709'
Set-TextCell 8 4 '0.2883'
Set-TextCell 8 5 '  -1.56%  '
Set-TextCell 9 4 '0.06458'
Set-TextCell 9 5 '  -1.09%  '
Set-TextCell 10 4 '21.98'
Set-TextCell 10 5 '  -0.34%  '
Set-TextCell 11 4 '0.07782'
Set-TextCell 11 5 '  +0.58%  '
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.871.77'
Set-TextCell 12 5 '  -0.33%  '
$ws.Cells.Item(13, 2).Value = 'Litecoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 13 4 '96.19'
Set-TextCell 13 5 '  -0.67%  '
Set-TextCell 14 4 '0.7243'
Set-TextCell 14 5 '  -2.35%  '
Set-TextCell 15 4 '5.136'
Set-TextCell 15 5 '  -1.33%  '
Set-TextCell 16 4 '280.84'
Set-TextCell 16 5 '  +2.29%  '
Set-TextCell 17 4 '30.366.22'
Set-TextCell 17 5 '  -1.14%  '
Set-TextCell 18 4 '13.02'
Set-TextCell 18 5 '  -1.81%  '
Set-TextCell 19 5 '  +0.01%  '
Set-TextCell 20 4 '0.000007496'
Set-TextCell 20 5 '  -0.59%  '
Set-TextCell 21 4 '2.114.18'
Set-TextCell 21 5 '  -0.36%  '
Set-TextCell 22 5 '  +0.05%  '
Set-TextCell 23 4 '5.244'
Set-TextCell 23 5 '  -0.44%  '
Set-TextCell 24 4 '6.231'
Set-TextCell 24 5 '  +0.35%  '
Set-TextCell 25 4 '163.75'
Set-TextCell 25 5 '  -0.93%  '
Set-TextCell 26 4 '9.055'
Set-TextCell 26 5 '  -1.64%  '
Set-TextCell 27 4 '18.72'
Set-TextCell 28 5 '  -1.96%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 29 4 '1.321'
Set-TextCell 29 5 '  -1.35%  '
$ws.Cells.Item(30, 2).Value = 'Stellar'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 30 4 '0.09621'
Set-TextCell 30 5 '  -2.50%  '
Set-TextCell 31 4 '1.486'
Set-TextCell 31 5 '  -1.21%  '
Set-TextCell 32 4 '4.228'
Set-TextCell 32 5 '  -1.48%  '
Set-TextCell 33 4 '4.117'
Set-TextCell 33 5 '  -0.05%  '
Set-TextCell 34 4 '0.04819'
Set-TextCell 34 5 '  -0.22%  '
Set-TextCell 35 5 '  -0.54%  '
Set-TextCell 36 4 '0.6906'
Set-TextCell 36 5 '  -0.97%  '
Set-TextCell 37 4 '2.712'
Set-TextCell 37 5 '  -0.17%  '
Set-TextCell 38 4 '0.01889'
Set-TextCell 38 5 '  +0.51%  '
Set-TextCell 39 4 '2.814'
Set-TextCell 39 5 '  +1.93%  '
Set-TextCell 40 4 '6.229'
Set-TextCell 40 5 '  -0.91%  '
Set-TextCell 41 4 '74.53'
Set-TextCell 41 5 '  +1.32%  '
Set-TextCell 42 4 '0.4223'
Set-TextCell 42 5 '  -0.56%  '
Set-TextCell 43 4 '1.929'
Set-TextCell 43 5 '  -3.16%  '
Set-TextCell 44 5 '  -0.11%  '
Set-TextCell 45 4 '0.8268'
Set-TextCell 45 5 '  -1.36%  '
Set-TextCell 46 4 '100.98'
Set-TextCell 46 5 '  -1.21%  '
Set-TextCell 47 4 '9.641'
Set-TextCell 47 5 '  +2.88%  '
Set-TextCell 48 4 '35.28'
Set-TextCell 48 5 '  -0.41%  '
Set-TextCell 49 4 '6.957'
Set-TextCell 49 5 '  -1.70%  '
Set-TextCell 50 4 '899.23'
Set-TextCell 50 5 '  -1.59%  '
Set-TextCell 51 4 '0.05725'
Set-TextCell 51 5 '  +0.47%  '
